# Auto-update predictions and index for 2025-10-26
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9 to make room for the Trabzonspor fixture.
# This shifts the former rows 9 (Inter Miami CF) and 10 (Struga Trim & Lum)
# down to rows 10 and 11 respectively, along with their formatting.
$ws.Rows.Item(9).Insert()

# ---- Row 2: Borussia Mönchengladbach - Bayern Munich ----
$ws.Range("A2").Value = "Borussia Mönchengladbach - Bayern Munich ✓: 0:3"
$ws.Range("B2").Value = "Bayern Munich"
$ws.Range("C2").Value = 73
$ws.Range("D2").ClearContents()
$ws.Range("E2").Value = 98
$ws.Range("F2").Value = 1.28
$ws.Range("G2").Value = "✓"

# ---- Row 3: Raja Club Athletic - Olympique Dcheira ----
$ws.Range("A3").Value = "Raja Club Athletic ✓ - Olympique Dcheira: 1:0"
$ws.Range("B3").Value = "Raja Club Athletic"
$ws.Range("C3").Value = 70
$ws.Range("D3").Value = 87
$ws.Range("E3").Value = 89
$ws.Range("F3").Value = 1.67
$ws.Range("G3").Value = "✓"

# ---- Row 4: Chelsea FC - Sunderland AFC ----
$ws.Range("A4").Value = "Chelsea FC X - Sunderland AFC: 1:2"
$ws.Range("B4").Value = "Chelsea FC"
$ws.Range("C4").Value = 64
$ws.Range("D4").ClearContents()
$ws.Range("E4").Value = 77
$ws.Range("F4").Value = 1.45
$ws.Range("G4").Value = "X"

# ---- Row 5: Inter Club d'Escaldes - FC Ordino ----
$ws.Range("A5").Value = "Inter Club d'Escaldes ✓ - FC Ordino: 2:1"
$ws.Range("B5").Value = "Inter Club d'Escaldes"
$ws.Range("C5").Value = 64
$ws.Range("D5").Value = 59
$ws.Range("E5").ClearContents()
$ws.Range("F5").Value = 2.5
$ws.Range("G5").Value = "✓"

# ---- Row 6: SL Benfica - FC Arouca ----
$ws.Range("A6").Value = "SL Benfica ✓ - FC Arouca: 5:0"
$ws.Range("B6").Value = "SL Benfica"
$ws.Range("C6").Value = 62
$ws.Range("D6").Value = 100
$ws.Range("E6").ClearContents()
$ws.Range("F6").Value = 1.18
$ws.Range("G6").Value = "✓"

# ---- Row 7: Levski Sofia - Dobrudzha Dobrich ----
$ws.Range("A7").Value = "Levski Sofia ✓ - Dobrudzha Dobrich: 3:0"
$ws.Range("B7").Value = "Levski Sofia"
$ws.Range("C7").Value = 59
$ws.Range("D7").Value = 78
$ws.Range("E7").ClearContents()
$ws.Range("F7").Value = 1.18
$ws.Range("G7").Value = "✓"

# ---- Row 8: Borussia Dortmund - 1.FC Köln ----
$ws.Range("A8").Value = "Borussia Dortmund ✓ - 1.FC Köln: 1:0"
$ws.Range("B8").Value = "Borussia Dortmund"
$ws.Range("C8").Value = 57
$ws.Range("D8").Value = 96
$ws.Range("E8").Value = 100
$ws.Range("F8").Value = 1.42
$ws.Range("G8").Value = "✓"

# ---- Row 9 (new): Trabzonspor - Eyüpspor ----
$ws.Range("A9").Value = "Trabzonspor ✓ - Eyüpspor: 2:0"
$ws.Range("B9").Value = "Trabzonspor"
$ws.Range("C9").Value = 56
$ws.Range("D9").Value = 100
$ws.Range("E9").ClearContents()
$ws.Range("F9").Value = 1.62
$ws.Range("G9").Value = "✓"

# ---- Row 10 (was row 9): Inter Miami CF - Nashville SC (unchanged content) ----
$ws.Range("A10").Value = "Inter Miami CF ✓ - Nashville SC: 3:1"
$ws.Range("B10").Value = "Inter Miami CF"
$ws.Range("C10").Value = 56
$ws.Range("D10").Value = 59
$ws.Range("E10").ClearContents()
$ws.Range("F10").Value = 2.5
$ws.Range("G10").Value = "✓"

# ---- Row 11 (was row 10): Struga Trim & Lum - AP Brera Strumica ----
$ws.Range("A11").Value = "Struga Trim & Lum X - AP Brera Strumica: 2:4"
$ws.Range("B11").Value = "Struga Trim & Lum"
$ws.Range("C11").Value = 53
$ws.Range("D11").ClearContents()
$ws.Range("E11").Value = 94
$ws.Range("F11").Value = 1.83
$ws.Range("G11").Value = "X"
